$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format A2:E7 as Text so the numeric-looking strings are stored as text,
# matching the inlineStr cells in the target workbook (not converted to numbers).
$ws.Range("A2:E7").NumberFormat = "@"

# Rows 2-7: new benchmark timing values (columns A-E), replacing the old readings.
$ws.Range("A2").Value = "62.33076588000003"
$ws.Range("B2").Value = "78.92645529999989"
$ws.Range("C2").Value = "25.441004659999855"
$ws.Range("D2").Value = "38.781792619999855"
$ws.Range("E2").Value = "41.52080139999992"

$ws.Range("A3").Value = "122.09017631999977"
$ws.Range("B3").Value = "166.12965471999976"
$ws.Range("C3").Value = "50.49119200000014"
$ws.Range("D3").Value = "79.59224567999911"
$ws.Range("E3").Value = "83.35338906000061"

$ws.Range("A4").Value = "249.80314116000002"
$ws.Range("B4").Value = "351.9087772999988"
$ws.Range("C4").Value = "106.55403948000043"
$ws.Range("D4").Value = "173.29918742000075"
$ws.Range("E4").Value = "169.911224439999"

$ws.Range("A5").Value = "509.63340572000135"
$ws.Range("B5").Value = "723.9605081800008"
$ws.Range("C5").Value = "216.77586504000146"
$ws.Range("D5").Value = "353.40306548000285"
$ws.Range("E5").Value = "347.88400354000146"

$ws.Range("A6").Value = "1069.528617359997"
$ws.Range("B6").Value = "1493.349852399998"
$ws.Range("C6").Value = "459.8563705799995"
$ws.Range("D6").Value = "757.420346040002"
$ws.Range("E6").Value = "742.149670959999"

$ws.Range("A7").Value = "2118.5710676199974"
$ws.Range("B7").Value = "2970.9002986600035"
$ws.Range("C7").Value = "908.8883996800018"
$ws.Range("D7").Value = "1488.6360109800057"
$ws.Range("E7").Value = "1461.5822356599938"

# Row 17, column F: n value doubled from 65536 to 131072.
$ws.Range("F17").Value = 131072

# New rows 18-22: extra "n" values extending the results table.
$ws.Range("F18").Value = 262144
$ws.Range("F19").Value = 524288
$ws.Range("F20").Value = 1048576
$ws.Range("F21").Value = 2097152
$ws.Range("F22").Value = 4194304
